$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D8").Value = 12196500
$ws.Range("E8").Value = 5819900
$ws.Range("F8").Value = 10491500
$ws.Range("G8").Value = 5068500
$ws.Range("H8").Value = 1467400
$ws.Range("I8").Value = 1285600
$ws.Range("J8").Value = 1232700
$ws.Range("D15").Value = -140000
$ws.Range("E15").Value = -69500
$ws.Range("F15").Value = -128400
$ws.Range("G15").Value = -62100
$ws.Range("D17").Value = 6936600
$ws.Range("E17").Value = 3382300
$ws.Range("F17").Value = 5953300
$ws.Range("G17").Value = 2907100
$ws.Range("H17").Value = 814200
$ws.Range("I17").Value = 768400
$ws.Range("J17").Value = 735400
$ws.Range("D18").Value = 5259900
$ws.Range("E18").Value = 2437700
$ws.Range("F18").Value = 4538200
$ws.Range("G18").Value = 2161400
$ws.Range("H18").Value = 653300
$ws.Range("I18").Value = 517200
$ws.Range("J18").Value = 497300
$ws.Range("D20").Value = -1252900
$ws.Range("E20").Value = -638000
$ws.Range("F20").Value = -1357500
$ws.Range("G20").Value = -746500
$ws.Range("H20").Value = -348300
$ws.Range("I20").Value = -124400
$ws.Range("J20").Value = -167800
$ws.Range("D21").Value = 4146900
$ws.Range("E21").Value = 1869100
$ws.Range("F21").Value = 3309100
$ws.Range("G21").Value = 1477000
$ws.Range("D23").Value = 4007000
$ws.Range("E23").Value = 1799600
$ws.Range("F23").Value = 3180700
$ws.Range("G23").Value = 1414900
$ws.Range("H23").Value = 305000
$ws.Range("I23").Value = 392800
$ws.Range("J23").Value = 329500
$ws.Range("D24").Value = 1421000
$ws.Range("E24").Value = 648600
$ws.Range("F24").Value = 1145600
$ws.Range("G24").Value = 519600
$ws.Range("H24").Value = 105800
$ws.Range("I24").Value = 124000
$ws.Range("J24").Value = 103900
$ws.Range("D26").Value = 2585900
$ws.Range("E26").Value = 1151000
$ws.Range("F26").Value = 2035100
$ws.Range("G26").Value = 895300
$ws.Range("H26").Value = 199100
$ws.Range("I26").Value = 268800
$ws.Range("J26").Value = 225600
$ws.Range("D27").Value = 2581300
$ws.Range("E27").Value = 1149100
$ws.Range("F27").Value = 2032100
$ws.Range("G27").Value = 893700
$ws.Range("H27").Value = 194600
$ws.Range("I27").Value = 268800
$ws.Range("J27").Value = 225600
$ws.Range("D32").Value = 1252900
$ws.Range("E32").Value = 638000
$ws.Range("F32").Value = 1357500
$ws.Range("G32").Value = 746500
$ws.Range("H32").Value = 348300
$ws.Range("I32").Value = 124400
$ws.Range("J32").Value = 167800
$ws.Range("D33").Value = 2581300
$ws.Range("E33").Value = 1149100
$ws.Range("F33").Value = 2032100
$ws.Range("G33").Value = 893700
$ws.Range("H33").Value = 194600
$ws.Range("I33").Value = 268800
$ws.Range("J33").Value = 225600
$ws.Range("D35").Value = 2581300
$ws.Range("E35").Value = 1149100
$ws.Range("F35").Value = 2032100
$ws.Range("G35").Value = 893700
$ws.Range("H35").Value = 194600
$ws.Range("I35").Value = 268800
$ws.Range("J35").Value = 225600
$ws.Range("D41").Value = 8302200
$ws.Range("E41").Value = 6594700
$ws.Range("H41").Value = 3944700
$ws.Range("I41").Value = 3574500
$ws.Range("J41").Value = 3862600
$ws.Range("D42").Value = 9399300
$ws.Range("D48").Value = 563500
$ws.Range("E48").Value = 551600
$ws.Range("H48").Value = 390900
$ws.Range("I48").Value = 371500
$ws.Range("J48").Value = 361500
$ws.Range("D49").Value = 1083600
$ws.Range("E49").Value = 1083600
$ws.Range("D54").Value = 164371300
$ws.Range("E54").Value = 144208100
$ws.Range("H54").Value = 57888000
$ws.Range("I54").Value = 55487200
$ws.Range("J54").Value = 54568400
$ws.Range("D57").Value = 6607600
$ws.Range("E57").Value = 6715600
$ws.Range("D61").Value = 13489800
$ws.Range("E61").Value = 13317200
$ws.Range("H61").Value = 4772800
$ws.Range("I61").Value = 4567100
$ws.Range("J61").Value = 4478600
$ws.Range("D66").Value = 147392900
$ws.Range("E66").Value = 128343500
$ws.Range("H66").Value = 52651400
$ws.Range("I66").Value = 50363100
$ws.Range("J66").Value = 49746700
$ws.Range("D72").Value = 10067000
$ws.Range("E72").Value = 8640600
$ws.Range("H72").Value = 5167800
$ws.Range("I72").Value = 5055600
$ws.Range("J72").Value = 4753300
$ws.Range("D76").Value = 16978400
$ws.Range("E76").Value = 15864600
$ws.Range("H76").Value = 5236600
$ws.Range("I76").Value = 5124100
$ws.Range("J76").Value = 4821700
$ws.Range("D81").Value = 2581300
$ws.Range("E81").Value = 1149100
$ws.Range("F81").Value = 2032100
$ws.Range("G81").Value = 893700
$ws.Range("H81").Value = 194600
$ws.Range("I81").Value = 268800
$ws.Range("J81").Value = 225600
$ws.Range("D83").Value = 140000
$ws.Range("E83").Value = 69500
$ws.Range("F83").Value = 128400
$ws.Range("G83").Value = 62100
$ws.Range("D89").Value = 1363200
$ws.Range("E89").Value = -600400
$ws.Range("F89").Value = 5484300
$ws.Range("G89").Value = 2250700
$ws.Range("D91").Value = -132800
$ws.Range("E91").Value = -61700
$ws.Range("F91").Value = -182600
$ws.Range("G91").Value = -97700
$ws.Range("D94").Value = -29618800
$ws.Range("E94").Value = -11319200
$ws.Range("F94").Value = -22820800
$ws.Range("G94").Value = -11543200
$ws.Range("D96").Value = -498700
$ws.Range("E96").Value = -492900
$ws.Range("F96").Value = -423400
$ws.Range("G96").Value = -419200
$ws.Range("D100").Value = 30322800
$ws.Range("E100").Value = 12292200
$ws.Range("F100").Value = 18132500
$ws.Range("G100").Value = 9478700
$ws.Range("D101").Value = 6900
$ws.Range("E101").Value = -5900
$ws.Range("F101").Value = -29000
$ws.Range("G101").Value = 5100
$ws.Range("D102").Value = 2074200
$ws.Range("E102").Value = 366700
$ws.Range("F102").Value = 766900
$ws.Range("G102").Value = 191400
